$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MEJORAR")

# Add the new product code as a new row at the end of the list (row 44),
# matching the style/format of the other regular data rows (e.g. A2).
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A44").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0
$ws.Range("A44").Value = "EVOL3975"

# Update the active selection to reflect the new end-of-list cell, as in the diff.
$ws.Activate()
$ws.Range("A45").Select() | Out-Null
